$d = $word.ActiveDocument

# 1) Remove the leading "All page/line/reference numbers refer to the tracked
#    revised manuscript." paragraph and the blank paragraph that follows it.
$p1 = $d.Paragraphs.Item(1)
$p2 = $d.Paragraphs.Item(2)
$removeRange = $d.Range($p1.Range.Start, $p2.Range.End)
$removeRange.Delete()

# 2) Change "As per instructions, ..." to "As per previous instructions, ..."
#    by inserting the word "previous " right after "As per ". Toggling a
#    character property on the newly inserted text and then reverting it
#    forces Word to keep the three pieces of text ("As per ", "previous ",
#    and "instructions, ...") as separate runs instead of silently
#    re-merging them into a single run.
$findRange = $d.Content
$found = $findRange.Find.Execute("As per ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$insertionPoint = $d.Range($findRange.End, $findRange.End)
$insertionPoint.InsertBefore("previous ")

$insertedRange = $d.Range($findRange.End, $findRange.End + 9)
$insertedRange.Font.Bold = 1
$insertedRange.Font.Bold = 0
